$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking
# strings (e.g. "1.012") are stored as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.305.14'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').Value = '1.937.95'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').Value = '321.40'
$ws.Range('E5').Value = '  -1.93%  '
$ws.Range('D6').Value = '1.011'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').Value = '0.4757'
$ws.Range('E7').Value = '  -4.45%  '
$ws.Range('D8').Value = '0.4066'
$ws.Range('E8').Value = '  -3.58%  '
$ws.Range('D9').Value = '53.54'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('D10').Value = '0.08528'
$ws.Range('E10').Value = '  -7.59%  '
$ws.Range('D11').Value = '1.053'
$ws.Range('E11').Value = '  -4.10%  '
$ws.Range('D12').Value = '22.34'
$ws.Range('E12').Value = '  -2.39%  '
$ws.Range('D13').Value = '1.937.44'
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('D14').Value = '7.537'
$ws.Range('E14').Value = '  -3.98%  '
$ws.Range('D15').Value = '6.132'
$ws.Range('E15').Value = '  -4.90%  '
$ws.Range('D16').Value = '1.013'
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('D17').Value = '90.03'
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('D18').Value = '0.00001071'
$ws.Range('E18').Value = '  -2.68%  '
$ws.Range('D19').Value = '0.06613'
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').Value = '18.26'
$ws.Range('E20').Value = '  -5.20%  '
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').Value = '5.815'
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('D23').Value = '28.351.43'
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('D24').Value = '11.46'
$ws.Range('E24').Value = '  -4.96%  '
$ws.Range('D25').Value = '2.298'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').Value = '2.225.97'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').Value = '155.41'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').Value = '20.24'
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('D29').Value = '2.174'
$ws.Range('E29').Value = '  -3.83%  '
$ws.Range('D30').Value = '5.775'
$ws.Range('E30').Value = '  -8.02%  '
$ws.Range('D31').Value = '123.98'
$ws.Range('E31').Value = '  -1.84%  '
$ws.Range('D32').Value = '0.9870'
$ws.Range('E32').Value = '  -5.58%  '
$ws.Range('D33').Value = '0.09615'
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.450'
$ws.Range('E34').Value = '  -5.04%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '3.673'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').Value = '5.606'
$ws.Range('E36').Value = '  -3.41%  '
$ws.Range('D37').Value = '9.288'
$ws.Range('E37').Value = '  +2.94%  '
$ws.Range('D38').Value = '0.02325'
$ws.Range('E38').Value = '  -4.13%  '
$ws.Range('D39').Value = '0.06176'
$ws.Range('E39').Value = '  -2.82%  '
$ws.Range('D40').Value = '1.242'
$ws.Range('E40').Value = '  -6.28%  '
$ws.Range('D41').Value = '0.6217'
$ws.Range('E41').Value = '  -3.51%  '
$ws.Range('D42').Value = '11.19'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').Value = '0.1915'
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('D45').Value = '1.330'
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('D46').Value = '0.5936'
$ws.Range('E46').Value = '  -4.87%  '
$ws.Range('D47').Value = '12.90'
$ws.Range('E47').Value = '  -2.96%  '
$ws.Range('D48').Value = '2.053'
$ws.Range('E48').Value = '  -6.58%  '
$ws.Range('D49').Value = '3.398'
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('D50').Value = '0.06792'
$ws.Range('E50').Value = '  -2.87%  '
$ws.Range('D51').Value = '110.15'
$ws.Range('E51').Value = '  -1.86%  '

# Restore the default cell style (removes the temporary text-format
# styling so the cells keep their original look/style index).
$ws.Range("D2:E51").Style = "Normal"
